# Apply the commit's changes to the workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Metadata sheet: bump the "Date" row value.
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# ---------------------------------------------------------------------
# 2) Elements sheet: swap the two "Mapping" columns (AK <-> AL), i.e.
#    the "Mapping: RIM Mapping" column and the
#    "Mapping: Spécification métier vers l'extension ROR
#    AvailableTimeNumberDaysofWeek" column change places - the French
#    business-mapping column now comes first (AK), RIM Mapping second (AL).
# ---------------------------------------------------------------------
$wsElem = $wb.Worksheets.Item("Elements")

$lastRow = $wsElem.Cells.Item($wsElem.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $wsElem.Cells.Item($r, 37)
    $alCell = $wsElem.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    if ($akVal -ne $alVal) {
        if ($alVal -eq "") {
            $akCell.ClearContents()
        } else {
            $akCell.Value = $alVal
        }
        if ($akVal -eq "") {
            $alCell.ClearContents()
        } else {
            $alCell.Value = $akVal
        }
    }
}

# Column widths follow their new contents: the wide French mapping text
# now lives in column AK (37, ~90.64 chars wide), and the narrower
# "RIM Mapping" text in AL (38, ~24.98 chars wide).
$wsElem.Columns.Item(37).ColumnWidth = 89.8
$wsElem.Columns.Item(38).ColumnWidth = 24.15
